# "Generate Report for Handoff"
# Updates the localization-status report: moves rows from "In Translation" to
# "Ready for handoff" and refreshes the handoff timestamps, then widens the
# "Status" columns so the new (longer) status text fits.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps ---
$wsOverview.Range("G2").Value = "2017-02-21 05:05:06"
$wsDeDe.Range("H2").Value = "2017-02-21 05:05:06"
$wsZhCn.Range("H2").Value = "2017-02-21 05:04:51"

# --- Widen the Status columns to fit "Ready for handoff" ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
